$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProviderOptions")

# Insert a new column before column R, shifting TASK_TAB..SCROLL_INTO_ELEMENT
# (and all data below them) one column to the right.
$ws.Columns("R").Insert()

# Populate the newly inserted column R with the new "VIEWALL" field.
$ws.Range("R1").Value = "VIEWALL"
$ws.Range("R7").Value = "Click"
$ws.Range("R8").Value = "Click"
$ws.Range("R9").Value = "Click"
$ws.Range("R10").Value = "Click"

# Update the view state to match where the user ended up after editing
# (scroll so column N is visible, then land the selection on R11).
$excel.Goto($ws.Range("N1"), $true)
$ws.Range("R11").Select()
